$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = -0.5990164247158317
    "C2" = 0.6519872051416672
    "D2" = -0.7361883853993072
    "B3" = -0.7088104133122732
    "C3" = 0.7076565324976303
    "D3" = 0.6326077162823311
    "B4" = 0.744770647955053
    "C4" = -0.5465145445234002
    "D4" = -0.8296348842292615
    "B5" = 0.6801305272088517
    "C5" = -0.7344895769520807
    "D5" = 0.8250443118658586
    "B6" = 0.6440322597549629
    "C6" = 0.6146163731836789
    "D6" = -0.6407777576641157
    "B7" = 0.6521810734325758
    "C7" = 0.6913387624508864
    "D7" = 0.7207719602344737
    "B8" = -0.6793870527162835
    "C8" = -0.6686907473491526
    "D8" = 0.6454434791267544
    "B9" = -0.7349716221433342
    "C9" = 0.6061008255280288
    "D9" = 0.6230625624320384
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
